$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-parsed as a number by Excel;
# force Text format first so the stored value stays an exact string match.
$textCells = @("D5", "D6", "D11", "D14", "D19", "D20", "D21", "D22", "D25", "D27", "D31", "D32", "D35", "D36", "D38", "D39", "D43", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.440.69'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '2.649.28'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '596.96'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').Value = '158.94'
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('D9').Value = '2.648.16'
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('E10').Value = '  -1.93%  '
$ws.Range('D11').Value = '0.157'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').Value = '28.16'
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('D15').Value = '3.133.25'
$ws.Range('E16').Value = '  -2.40%  '
$ws.Range('D17').Value = '68.340.14'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '2.647.81'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').Value = '11.62'
$ws.Range('E19').Value = '  +2.20%  '
$ws.Range('D20').Value = '364.54'
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').Value = '7.48'
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('D22').Value = '4.42'
$ws.Range('E22').Value = '  +1.49%  '
$ws.Range('E23').Value = '  -1.21%  '
$ws.Range('E24').Value = '  +1.21%  '
$ws.Range('D25').Value = '74.99'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '9.98'
$ws.Range('E27').Value = '  +1.96%  '
$ws.Range('E29').Value = '  -2.18%  '
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').Value = '570.27'
$ws.Range('E31').Value = '  +1.51%  '
$ws.Range('D32').Value = '8.06'
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('D35').Value = '1.64'
$ws.Range('E35').Value = '  +4.02%  '
$ws.Range('D36').Value = '0.129'
$ws.Range('E36').Value = '  -1.17%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').Value = '160.85'
$ws.Range('E38').Value = '  -0.63%  '
$ws.Range('D39').Value = '19.66'
$ws.Range('E39').Value = '  +1.81%  '
$ws.Range('E40').Value = '  -0.82%  '
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '2.65'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('E44').Value = '  -5.11%  '
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').Value = '158.56'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D47').Value = '3.82'
$ws.Range('E47').Value = '  +1.86%  '
$ws.Range('D48').Value = '21.88'
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('D51').Value = '0.576'
$ws.Range('E51').Value = '  +2.30%  '
